# INSERIDO EXAME VULVOSCOPIA MOARCY
#
# A new exam row is inserted right above the existing "COLPOSCOPIA (Dr. MOARCY)"
# row (row 479) for the "Instituto de Diagn\u00f3stico da Mulher" location, pushing
# COLPOSCOPIA and every row below it down by one. The new row repeats the same
# Local and price (90) as the COLPOSCOPIA row that used to sit at that position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A real full-row insert pushes every row below it down by one, including the
# (empty) row sitting at the very bottom of the worksheet, which simply falls
# off the grid. Drop that trailing empty row first so the shift below mirrors
# what Excel itself does.
$ws.Rows.Item(1048576).Delete()

# Insert a new blank row at row 479, shifting existing row 479 (and below) down to 480.
$ws.Rows.Item(479).Insert()

$local = "Instituto de Diagn" + [char]0xF3 + "stico da Mulher"

$ws.Range("A479").Value = $local
$ws.Range("B479").Value = "VULVOSCOPIA (Dr. MOARCY)"
$ws.Range("C479").Value = 90

# Match the row height used by its neighbours (13.8pt) instead of the sheet's
# default (14.25pt).
$ws.Rows.Item(479).RowHeight = 13.8
